$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$rowValues = @(1.02, 1.04777558345717, 1.05671684321732, 1.045195325001565, 1.067280212569553, 1.050095241103683, 1.052822855969226, 1.059453789623616, 1.047964246794331, 1.06998849922683, 1.021365751152969)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(2, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.048830217025427, 1.057590295161704, 1.046097723799521, 1.068313149488422, 1.050460131422822, 1.053525582428084, 1.060140783526407, 1.048677799829911, 1.070836650192822, 1.021607595383811)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(3, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.049512764553604, 1.05815558199185, 1.046682038034412, 1.068981990817296, 1.050695045803165, 1.053979831083095, 1.060584779143082, 1.049139285749169, 1.071385312456849, 1.021763742538002)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(4, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.049799738520976, 1.058393253217945, 1.046927779433316, 1.069263281921324, 1.05079351778988, 1.05417068628231, 1.060771306424306, 1.049333238882636, 1.07161593426001, 1.021829304670868)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(5, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.049847924532424, 1.058433160693217, 1.046969046092176, 1.069310518386261, 1.050810034892885, 1.054202725197636, 1.060802617632844, 1.049365801200076, 1.071654654565323, 1.021840308019508)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(6, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.049516598991926, 1.058158757670196, 1.046685321267636, 1.068985749010255, 1.050696362714723, 1.05398238173581, 1.060587272034943, 1.049141877579112, 1.071388394177822, 1.021764618905394)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(7, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.048131975129475, 1.057012008347463, 1.045500211381033, 1.067629202514624, 1.050218804636236, 1.053060441250762, 1.059686073098632, 1.048205442928615, 1.070275166275649, 1.021447554285295)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(8, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.045693075589593, 1.054992115514987, 1.043414996615316, 1.065242351024981, 1.049368148127516, 1.051432333870454, 1.058093958316345, 1.046553573948524, 1.068312397216843, 1.020886233510831)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(9, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.044067796017854, 1.053646102804613, 1.042026966104038, 1.063653529656857, 1.048794907777611, 1.050344567606856, 1.057029816428239, 1.045451168346709, 1.067003150273303, 1.020510274334558)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(10, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.043364184192717, 1.053063407652107, 1.041426439961358, 1.062966128538798, 1.048545234234015, 1.049872994668029, 1.056568385958126, 1.044973542688881, 1.066436060873909, 1.020347068022708)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(11, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.043102852857599, 1.052846989669472, 1.041203453316935, 1.062710882853301, 1.048452275542818, 1.049697746971489, 1.056396892613601, 1.044796089851824, 1.066225392214782, 1.020286383948305)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(12, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.043158908352185, 1.052893411108544, 1.041251281286798, 1.062765630038004, 1.048472225399201, 1.049735342020422, 1.05643368292909, 1.044834155958376, 1.06627058254965, 1.020299403690866)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(13, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.043342582027992, 1.053045518045074, 1.041408006273495, 1.062945028101424, 1.048537554703284, 1.049858510362338, 1.056554212253421, 1.044958875225552, 1.066418647474575, 1.020342053121809)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(14, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.043455752270213, 1.053139238933489, 1.041504579725793, 1.06305557253102, 1.048577777248811, 1.049934387246546, 1.056628461429447, 1.045035713401584, 1.066509871667304, 1.02036832263466)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(15, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.019999999999999, 1.044114495441213, 1.053684777281788, 1.042066831615407, 1.063699162211411, 1.048811447100211, 1.05037585249449, 1.057060426365277, 1.045482860972808, 1.067040782417069, 1.020521097099896)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(16, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.044527746186949, 1.054027016044474, 1.042419651543322, 1.064103021503173, 1.048957632048421, 1.050652621306975, 1.0573312126251, 1.045763270860348, 1.067373761473611, 1.020616817877662)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(17, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.044768802364759, 1.054226651343833, 1.042625493679326, 1.064338640683376, 1.049042758703875, 1.050814001549689, 1.057489094943477, 1.045926802488479, 1.067567965451037, 1.020672610280872)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(18, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.044850998678979, 1.054294724054256, 1.042695688700193, 1.064418990038556, 1.04907176085263, 1.050869018825799, 1.057542918090591, 1.045982558010888, 1.067634181087971, 1.020691627275347)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(19, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.044483406846256, 1.053990295694397, 1.042381792311682, 1.064059685544326, 1.048941962324595, 1.050622932251016, 1.057302166309078, 1.045733188295096, 1.067338037736787, 1.020606552070823)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(20, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.043288494119787, 1.053000725762739, 1.041361852584103, 1.062892197475811, 1.048518322891776, 1.049822242687863, 1.05651872206254, 1.044922149651357, 1.066375046772255, 1.020329495637312)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(21, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.042537328850169, 1.052378665192829, 1.040721013009401, 1.062158647306273, 1.048250697955976, 1.049318328465843, 1.056025575401216, 1.044411978019499, 1.065769423432875, 1.020154940714746)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(22, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.042935524155165, 1.052708419662497, 1.041060692652194, 1.062547469053389, 1.04839269099031, 1.049585509249039, 1.056287055208512, 1.044682452291142, 1.066090490243805, 1.020247509520577)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(23, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.044503441842266, 1.054006887999082, 1.042398899119718, 1.064079267029066, 1.048949043232509, 1.050636347625209, 1.057315291282733, 1.045746781393918, 1.067354179808409, 1.020611190870339)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(24, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.046323473791962, 1.055514205546683, 1.043953704310738, 1.065858985633248, 1.049589145329922, 1.051853655765146, 1.058506040859164, 1.046980827385799, 1.068819950824411, 1.021031656826332)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(25, $colIndexes[$i]).Value = $rowValues[$i]
}
